$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q1"), matching the sheet order:
#    总计, 2022-Q4, 2022-Q1, 2021-Q4, 2021-Q3
# ------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$oldQ1Sheet   = $wb.Worksheets.Item(2)

$q4Sheet = $wb.Worksheets.Add($oldQ1Sheet, $null)
$q4Sheet.Name = "2022-Q4"

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the 2022-Q4 row at the
#    top of the data and push the existing rows down by one.
# ------------------------------------------------------------------

# Row 5 (new) <- old row 4 data (2021-Q3 / 6 / 0.11)
$summarySheet.Cells.Item(5, 1).Value = 3
$summarySheet.Cells.Item(4, 1).Copy()
$summarySheet.Cells.Item(5, 1).PasteSpecial(-4122)
$summarySheet.Cells.Item(5, 2).Value = "2021-Q3"
$summarySheet.Cells.Item(5, 3).Value = 6
$summarySheet.Cells.Item(5, 4).Value = 0.11

# Row 4 <- old row 3 data (2021-Q4 / 4 / 0.78)
$summarySheet.Cells.Item(4, 2).Value = "2021-Q4"
$summarySheet.Cells.Item(4, 3).Value = 4
$summarySheet.Cells.Item(4, 4).Value = 0.78

# Row 3 <- old row 2 data (2022-Q1 / 3 / 0.28)
$summarySheet.Cells.Item(3, 2).Value = "2022-Q1"
$summarySheet.Cells.Item(3, 3).Value = 3
$summarySheet.Cells.Item(3, 4).Value = 0.28

# Row 2 <- new 2022-Q4 data (2022-Q4 / 1 / 0.09)
$summarySheet.Cells.Item(2, 2).Value = "2022-Q4"
$summarySheet.Cells.Item(2, 3).Value = 1
$summarySheet.Cells.Item(2, 4).Value = 0.09

# ------------------------------------------------------------------
# 3. Populate the new "2022-Q4" sheet with its fund-holding table,
#    copying header/index cell formatting from the sibling "2022-Q1"
#    sheet so the look matches the rest of the workbook.
# ------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4Sheet.Cells.Item(1, $col).Value = $headers[$col - 2]
    $oldQ1Sheet.Cells.Item(1, $col).Copy()
    $q4Sheet.Cells.Item(1, $col).PasteSpecial(-4122)
}

$q4Sheet.Cells.Item(2, 1).Value = 0
$oldQ1Sheet.Cells.Item(2, 1).Copy()
$q4Sheet.Cells.Item(2, 1).PasteSpecial(-4122)

$q4Sheet.Cells.Item(2, 2).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 2).Value = "001144"

$q4Sheet.Cells.Item(2, 3).Value = "大成互联网思维混合"

$q4Sheet.Cells.Item(2, 4).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 4).Value = "2.35"

$q4Sheet.Cells.Item(2, 5).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 5).Value = "68.92"

$q4Sheet.Cells.Item(2, 6).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 6).Value = "3.82"

$q4Sheet.Cells.Item(2, 7).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 7).Value = "0.0898"

$q4Sheet.Cells.Item(2, 8).Value = 8
